$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Clear the old "pricing table" block (M1:R11) - its content is being
#    replaced by a new "Display Pins" table (M:N) and the pricing table
#    shifted right to P:U, with an extra row inserted.
# ---------------------------------------------------------------------------
$ws.Range("M1:R11").ClearContents()

# ---------------------------------------------------------------------------
# 2) Slider pin renumbering (left-hand table, columns A-K unaffected in
#    position, just renumbered Output Pin / Pin values).
# ---------------------------------------------------------------------------
$ws.Range("D2").Value = 10
$ws.Range("D3").Value = 9
$ws.Range("D4").Value = 8
$ws.Range("D5").Value = 7
$ws.Range("D6").Value = 6
$ws.Range("D7").Value = 5

$ws.Range("I5").Value = 17
$ws.Range("I6").Value = 1
$ws.Range("I7").Value = 0

# ---------------------------------------------------------------------------
# 3) New "Display Pins" table in columns M:N (rows 1-8).
# ---------------------------------------------------------------------------
$ws.Range("M1").Value = "Display Pins"

$ws.Range("M2").Value = "VCC"
$ws.Range("N2").Value = "ICSP 5V"

$ws.Range("M3").Value = "GND"
$ws.Range("N3").Value = "ICSP GND"

$ws.Range("M4").Value = "SCL"
$ws.Range("N4").Value = "ICSP SCK"

$ws.Range("M5").Value = "SDA"
$ws.Range("N5").Value = "ICSP COPI"

$ws.Range("M6").Value = "RS/DC"
$ws.Range("N6").Value = 2

$ws.Range("M7").Value = "RES"
$ws.Range("N7").Value = 3

$ws.Range("M8").Value = "CS"
$ws.Range("N8").Value = 4

# ---------------------------------------------------------------------------
# 4) Pricing table, shifted to P:U, with a new row for the "Lever-T"
#    potentiometer inserted as row 5 (pushing Schiebeknopf/Knopf/
#    Widerstand/Display down by one row).
# ---------------------------------------------------------------------------
$ws.Range("P1").Value = "Komponente"
$ws.Range("Q1").Value = "Anzahl"
$ws.Range("R1").Value = "St" + [char]0xFC + "ck (" + [char]0x20AC + ")"
$ws.Range("S1").Value = "Preis (" + [char]0x20AC + ")"
$ws.Range("U1").Value = "EAN"

$ws.Range("P2").Value = "Arduino Micro"
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 11.9
$ws.Range("S2").Formula = "=PRODUCT(Q2, R2)"
$ws.Range("U2").Value = "C-4250236822907"

$ws.Range("P3").Value = "Kabel"
$ws.Range("Q3").Value = 1
$ws.Range("R3").Value = 0
$ws.Range("S3").Formula = "=PRODUCT(Q3, R3)"
$ws.Range("U3").Value = "Sollte ich haben"

$ws.Range("P4").Value = "Potentiometer"
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 7.49
$ws.Range("S4").Formula = "=PRODUCT(Q4, R4)"
$ws.Range("T4").Value = "Lever 1"
$ws.Range("U4").Value = "C-2050000114397"

$ws.Range("P5").Value = "Potentiometer"
$ws.Range("Q5").Value = 5
$ws.Range("R5").Value = 21.39
$ws.Range("T5").Value = "Lever-T"
$ws.Range("U5").Value = "C-2050000114243"

$ws.Range("P6").Value = "Schiebeknopf"
$ws.Range("Q6").Value = 0
$ws.Range("R6").Value = 1.49
$ws.Range("T6").Value = "Schwarz"
$ws.Range("U6").Value = "C-2050000113758"

$ws.Range("P7").Value = "Schiebeknopf"
$ws.Range("Q7").Value = 0
$ws.Range("R7").Value = 1.49
$ws.Range("T7").Value = "Blau"
$ws.Range("U7").Value = "C-2050000113840"

$ws.Range("P8").Value = "Knopf "
$ws.Range("Q8").Value = 0
$ws.Range("R8").Value = 2.29
$ws.Range("U8").Value = "C-2050004878899"

$ws.Range("P9").Value = "Widerstand "
$ws.Range("Q9").Value = 0
$ws.Range("R9").Value = 0.05
$ws.Range("U9").Value = "C-2050000096181"

$ws.Range("P10").Value = "Display"
$ws.Range("Q10").Value = 0
$ws.Range("R10").Value = 9.99
$ws.Range("U10").Value = "C-4250236809434"

# Fill the shared PRODUCT() formula down for the new row range S5:S10.
$ws.Range("S5:S10").Formula = "=PRODUCT(Q5, R5)"

$ws.Range("P12").Value = "Gesamt"
$ws.Range("S12").Formula = "=SUM(S2:S10)"

# ---------------------------------------------------------------------------
# 5) Styling. Reuse formatting already present elsewhere on the sheet by
#    copying it across (keeps fills/fonts/theme colors byte identical),
#    then nudge horizontal alignment to "left" for the new pin-label column.
# ---------------------------------------------------------------------------
function Copy-Style($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial(-4122) | Out-Null
}

# Headers (M1 matches the other black header cells, e.g. A1).
Copy-Style "A1" "M1"
Copy-Style "A1" "P1"
Copy-Style "A1" "Q1"
Copy-Style "A1" "R1"
Copy-Style "A1" "S1"
Copy-Style "A1" "U1"

# Display-pin table.
Copy-Style "B3" "M2"          # red fill, center/top
Copy-Style "B3" "N2"
$ws.Range("N2").HorizontalAlignment = -4131   # xlLeft

Copy-Style "A1" "M3"          # black fill / white font, center/top
Copy-Style "A1" "N3"
$ws.Range("N3").HorizontalAlignment = -4131   # xlLeft

Copy-Style "B4" "M4"          # purple fill, center/top
Copy-Style "B4" "N4"
$ws.Range("N4").HorizontalAlignment = -4131   # xlLeft

Copy-Style "B4" "M5"
Copy-Style "B4" "N5"
$ws.Range("N5").HorizontalAlignment = -4131   # xlLeft

Copy-Style "B4" "M6"
Copy-Style "B4" "N6"

Copy-Style "B4" "M7"
Copy-Style "B4" "N7"

Copy-Style "B4" "M8"
Copy-Style "B4" "N8"

# Pricing table (re-use the same styles the rows used before the shift).
Copy-Style "M5" "P2"   # themed accent fill, center/top
Copy-Style "N2" "Q2"
$ws.Range("Q2").HorizontalAlignment = -4108   # xlCenter (restore center; N2 donor was left-aligned)
$ws.Range("Q2").VerticalAlignment = -4160     # xlTop
Copy-Style "O2" "R2"
Copy-Style "P2" "S2"
Copy-Style "R2" "U2"

Copy-Style "M3" "P3"
Copy-Style "N3" "Q3"
Copy-Style "O3" "R3"
Copy-Style "P3" "S3"
Copy-Style "R2" "U3"

Copy-Style "B3" "P4"
Copy-Style "N2" "Q4"
$ws.Range("Q4").HorizontalAlignment = -4108
$ws.Range("Q4").VerticalAlignment = -4160
Copy-Style "O2" "R4"
Copy-Style "P2" "S4"
Copy-Style "R2" "T4"
Copy-Style "R2" "U4"

Copy-Style "A1" "P5"
Copy-Style "N2" "Q5"
$ws.Range("Q5").HorizontalAlignment = -4108
$ws.Range("Q5").VerticalAlignment = -4160
Copy-Style "O2" "R5"
Copy-Style "P2" "S5"
Copy-Style "R2" "T5"
Copy-Style "R2" "U5"

Copy-Style "M5" "P6"
Copy-Style "N2" "Q6"
$ws.Range("Q6").HorizontalAlignment = -4108
$ws.Range("Q6").VerticalAlignment = -4160
Copy-Style "O2" "R6"
Copy-Style "P2" "S6"
Copy-Style "R2" "T6"
Copy-Style "R2" "U6"

Copy-Style "M5" "P7"
Copy-Style "N2" "Q7"
$ws.Range("Q7").HorizontalAlignment = -4108
$ws.Range("Q7").VerticalAlignment = -4160
Copy-Style "O2" "R7"
Copy-Style "P2" "S7"
Copy-Style "R2" "T7"
Copy-Style "R2" "U7"

Copy-Style "M5" "P8"
Copy-Style "N2" "Q8"
$ws.Range("Q8").HorizontalAlignment = -4108
$ws.Range("Q8").VerticalAlignment = -4160
Copy-Style "O2" "R8"
Copy-Style "P2" "S8"
Copy-Style "R2" "U8"

Copy-Style "M5" "P9"
Copy-Style "N2" "Q9"
$ws.Range("Q9").HorizontalAlignment = -4108
$ws.Range("Q9").VerticalAlignment = -4160
Copy-Style "O2" "R9"
Copy-Style "P2" "S9"
Copy-Style "R2" "U9"

Copy-Style "M5" "P10"
Copy-Style "N2" "Q10"
$ws.Range("Q10").HorizontalAlignment = -4108
$ws.Range("Q10").VerticalAlignment = -4160
Copy-Style "O2" "R10"
Copy-Style "P2" "S10"
Copy-Style "R2" "U10"

Copy-Style "N2" "Q11"
$ws.Range("Q11").HorizontalAlignment = -4108
$ws.Range("Q11").VerticalAlignment = -4160

Copy-Style "P2" "P12"
Copy-Style "P2" "S12"

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 6) Page setup: landscape, paper size 9 (A4).
# ---------------------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 2
